# Update cryptocurrency price/volume data (and two rank swaps) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $text) {
    $cell = $ws.Range($cellAddress)
    # Prefix with an apostrophe so Excel stores numeric-looking strings
    # (e.g. "19.00", "0.0000150") as literal text instead of coercing them
    # to a Double and silently dropping significant trailing zeros.
    $cell.Value = "'" + $text
    # Re-apply the default "Normal" style so we do not leave a stray
    # quote-prefixed / text-formatted style behind on the cell.
    $cell.Style = "Normal"
}

Set-TextValue "D2" "56.781.52"
Set-TextValue "E2" "  -1.46%  "
Set-TextValue "D3" "2.992.06"
Set-TextValue "E3" "  -3.11%  "
Set-TextValue "E4" "  -0.23%  "
Set-TextValue "D5" "499.36"
Set-TextValue "E5" "  -3.31%  "
Set-TextValue "D6" "134.90"
Set-TextValue "E6" "  +4.28%  "
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "E8" "  -2.33%  "
Set-TextValue "D9" "7.24"
Set-TextValue "E9" "  +1.20%  "
Set-TextValue "E10" "  +1.07%  "
Set-TextValue "D11" "0.351"
Set-TextValue "E11" "  -3.98%  "
Set-TextValue "E12" "  -0.58%  "
Set-TextValue "D13" "3.501.38"
Set-TextValue "E13" "  -3.53%  "
Set-TextValue "D14" "25.21"
Set-TextValue "E14" "  +2.74%  "
Set-TextValue "D15" "56.742.06"
Set-TextValue "E15" "  -1.74%  "
Set-TextValue "D16" "0.0000150"
Set-TextValue "E16" "  +1.93%  "
Set-TextValue "D17" "2.987.51"
Set-TextValue "E17" "  -3.58%  "
Set-TextValue "D18" "5.68"
Set-TextValue "E18" "  +1.51%  "
Set-TextValue "E19" "  -2.67%  "
Set-TextValue "D20" "7.79"
Set-TextValue "E20" "  +1.28%  "
Set-TextValue "D21" "326.56"
Set-TextValue "E21" "  -2.51%  "
Set-TextValue "E22" "  -0.01%  "
Set-TextValue "D23" "0.469"
Set-TextValue "E23" "  -6.09%  "
Set-TextValue "D24" "62.08"
Set-TextValue "E24" "  -5.83%  "
Set-TextValue "D25" "0.999"
Set-TextValue "E25" "  -0.56%  "
Set-TextValue "E26" "  -1.76%  "
Set-TextValue "D27" "0.0₃0892"
Set-TextValue "E27" "  -0.18%  "
Set-TextValue "E28" "  -0.23%  "
Set-TextValue "E29" "  -3.74%  "
Set-TextValue "D30" "6.85"
Set-TextValue "E30" "  +2.35%  "
Set-TextValue "E31" "  -4.45%  "
Set-TextValue "B32" "Fetch.AI"
Set-TextValue "C32" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D32" "1.16"
Set-TextValue "E32" "  -5.90%  "
Set-TextValue "B33" "EthereumClassic"
Set-TextValue "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "20.31"
Set-TextValue "E33" "  -2.80%  "
Set-TextValue "D34" "155.66"
Set-TextValue "E34" "  -1.15%  "
Set-TextValue "D35" "4.44"
Set-TextValue "E35" "  -5.25%  "
Set-TextValue "E36" "  -3.96%  "
Set-TextValue "E37" "  -7.14%  "
Set-TextValue "E38" "  +1.26%  "
Set-TextValue "B39" "EnergySwap"
Set-TextValue "C39" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D39" "22.96"
Set-TextValue "E39" "  +1.05%  "
Set-TextValue "B40" "RenzoRestakedETH"
Set-TextValue "C40" "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue "D40" "3.023.82"
Set-TextValue "E40" "  -3.25%  "
Set-TextValue "D41" "36.41"
Set-TextValue "E41" "  -9.02%  "
Set-TextValue "E42" "  -0.37%  "
Set-TextValue "E43" "  -5.73%  "
Set-TextValue "D44" "2.233.01"
Set-TextValue "E44" "  -0.39%  "
Set-TextValue "D45" "0.995"
Set-TextValue "E45" "  -4.78%  "
Set-TextValue "E46" "  +0.41%  "
Set-TextValue "E47" "  -7.27%  "
Set-TextValue "E48" "  +12.16%  "
Set-TextValue "D49" "0.0236"
Set-TextValue "E49" "  +3.58%  "
Set-TextValue "E50" "  -5.19%  "
Set-TextValue "D51" "19.00"
Set-TextValue "E51" "  -4.82%  "
